$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Capture the two "NULL accession" rows (no accession in column A) before they move.
$row236 = @($ws1.Cells.Item(236,2).Value2, $ws1.Cells.Item(236,3).Value2, $ws1.Cells.Item(236,5).Value2, $ws1.Cells.Item(236,6).Value2, $ws1.Cells.Item(236,7).Value2, $ws1.Cells.Item(236,8).Value2, $ws1.Cells.Item(236,9).Value2, $ws1.Cells.Item(236,10).Value2)
$row237 = @($ws1.Cells.Item(237,2).Value2, $ws1.Cells.Item(237,3).Value2, $ws1.Cells.Item(237,5).Value2, $ws1.Cells.Item(237,6).Value2, $ws1.Cells.Item(237,7).Value2, $ws1.Cells.Item(237,8).Value2, $ws1.Cells.Item(237,9).Value2, $ws1.Cells.Item(237,10).Value2)

# Remove those two rows from Sheet1 - everything below shifts up by 2.
$ws1.Range("A236:J237").EntireRow.Delete()

# Create the new "NULL accessions" sheet right after Sheet1.
$wsNull = $wb.Worksheets.Add($null, $ws1)
$wsNull.Name = "NULL accessions"

$wsNull.Cells.Item(1,2).Value = $row236[0]
$wsNull.Cells.Item(1,3).Value = $row236[1]
$wsNull.Cells.Item(1,5).Value = $row236[2]
$wsNull.Cells.Item(1,6).Value = $row236[3]
$wsNull.Cells.Item(1,7).Value = $row236[4]
$wsNull.Cells.Item(1,8).Value = $row236[5]
$wsNull.Cells.Item(1,9).Value = $row236[6]
$wsNull.Cells.Item(1,10).Value = $row236[7]

$wsNull.Cells.Item(2,2).Value = $row237[0]
$wsNull.Cells.Item(2,3).Value = $row237[1]
$wsNull.Cells.Item(2,5).Value = $row237[2]
$wsNull.Cells.Item(2,6).Value = $row237[3]
$wsNull.Cells.Item(2,7).Value = $row237[4]
$wsNull.Cells.Item(2,8).Value = $row237[5]
$wsNull.Cells.Item(2,9).Value = $row237[6]
$wsNull.Cells.Item(2,10).Value = $row237[7]

# Sheet1 view tweaks (from the target workbook).
$ws1.Application.ActiveWindow.ScrollRow = 249
$ws1.Range("D276").Select()

# New sheet view selection.
$wsNull.Range("F27:G27").Select()

$ws1.Select()
